$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio2")

# --- Row 2: new language-probability split ---
$ws.Range("B2").Value = 0.07
$ws.Range("C2").Value = 0.2
$ws.Range("D2").Value = 0.25
$ws.Range("E2").Value = 0.25
$ws.Range("F2").Value = 0.15
$ws.Range("G2").Value = 0.08

# --- Row 4: new language-probability split ---
$ws.Range("B4").Value = 0.08
$ws.Range("C4").Value = 0.2
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 0.2
$ws.Range("F4").Value = 0.015
$ws.Range("G4").Value = 0.005

# --- Row 7: new language-probability split (same as row 4) ---
$ws.Range("B7").Value = 0.08
$ws.Range("C7").Value = 0.2
$ws.Range("D7").Value = 0.5
$ws.Range("E7").Value = 0.2
$ws.Range("F7").Value = 0.015
$ws.Range("G7").Value = 0.005

# G2/G4/G7 previously had no explicit number format (General); give them the
# same "0.000" style the rest of column G (and B:F) already carries.
$ws.Range("G2").NumberFormat = "0.000"
$ws.Range("G4").NumberFormat = "0.000"
$ws.Range("G7").NumberFormat = "0.000"

# --- Column I: running total check (should equal 1 for every language row) ---
$ws.Range("I2").Formula = "=SUM(B2:G2)"
$ws.Range("I3:I12").Formula = "=SUM(B3:G3)"
$ws.Range("I2:I12").NumberFormat = "0.000"

# Widen column I so the new total column is fully visible.
$ws.Columns.Item(9).AutoFit()

Write-Host "Job-offer language-probability generator columns populated."
